$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''328.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''-0.60%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''43.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''2.32%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.603'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''-1.51%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.08207'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''-1.69%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''8.768'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''-0.44%'
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''4.477'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''-1.03%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '''1.913'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''-5.73%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '''2.855'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-4.57%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").Value = '''0.9463'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''1.55%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.1220'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''-5.56%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = '''0.1922'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''-2.19%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '''0.09836'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''4.51%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.04459'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''15.36%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '''0.1070'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''0.61%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '''0.001282'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-1.24%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.006130'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''-0.88%'
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.486'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''1.18%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = '''8.803'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''5.65%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.1364'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''0.34%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = '''0.04415'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''0.05%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.001245'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''-0.53%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.004395'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''0.34%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.0001238'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''3.34%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.0004015'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''31.85%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = '''0.02774'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''-1.39%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.05743'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''3.06%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.007913'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''1.39%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.009940'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''11.47%'
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.1418'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''-1.77%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.002083'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-7.05%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.009746'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-13.77%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.00007278'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''3.87%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.00000000755'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''0.82%'
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.003382'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''6.51%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.002285'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''0.19%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.00002114'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''0.82%'
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.0002013'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''0.82%'
$ws.Range("E51").Style = "Normal"
